# Auto-generated edit script: refresh market-data-driven columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6735.5854
$ws.Range("I40").Value = 18398
$ws.Range("J40").Value = 4736.3145
$ws.Range("K40").Value = 18398
$ws.Range("L40").Value = 4736.3145
$ws.Range("M40").Value = -18223
$ws.Range("N40").Value = -5086.3145
$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3238
$ws.Range("H55").Value = 351.18182
$ws.Range("I55").Value = 175.16667
$ws.Range("J55").Value = 562.4
$ws.Range("K55").Value = 175.16667
$ws.Range("L55").Value = 562.4
$ws.Range("M55").Value = 38.83332999999999
$ws.Range("N55").Value = -990.4
$ws.Range("H60").Value = 1000
$ws.Range("J60").Value = 1000
$ws.Range("L60").Value = 3000
$ws.Range("N60").Value = -3968
$ws.Range("H70").Value = 4506.0557
$ws.Range("I70").Value = 838.1667
$ws.Range("K70").Value = 2514.5001
$ws.Range("M70").Value = -2244.5001
$ws.Range("H73").Value = 4506.0557
$ws.Range("I73").Value = 838.1667
$ws.Range("K73").Value = 2514.5001
$ws.Range("M73").Value = -1578.5001
$ws.Range("H135").Value = 1026.4736
$ws.Range("I135").Value = 861.2778
$ws.Range("K135").Value = 7751.500199999999
$ws.Range("M135").Value = -5216.500199999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9053972
$ws.Range("I45").Value = 15386993
$ws.Range("K45").Value = 15386993
$ws.Range("M45").Value = -15386616
$ws.Range("H61").Value = 9605.700000000001
$ws.Range("I61").Value = 11276.866
$ws.Range("J61").Value = 4592.2
$ws.Range("K61").Value = 11276.866
$ws.Range("L61").Value = 4592.2
$ws.Range("M61").Value = -11064.866
$ws.Range("N61").Value = -5016.2
$ws.Range("H136").Value = 9605.700000000001
$ws.Range("I136").Value = 11276.866
$ws.Range("J136").Value = 4592.2
$ws.Range("K136").Value = 33830.598
$ws.Range("L136").Value = 13776.6
$ws.Range("M136").Value = -31280.598
$ws.Range("N136").Value = -18876.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20836116
$ws.Range("I20").Value = 25642454
$ws.Range("J20").Value = 8649.666999999999
$ws.Range("K20").Value = 25642454
$ws.Range("L20").Value = 8649.666999999999
$ws.Range("M20").Value = -25642207
$ws.Range("N20").Value = -9143.666999999999
$ws.Range("H59").Value = 130000
$ws.Range("J59").Value = 130000
$ws.Range("L59").Value = 130000
$ws.Range("N59").Value = -131694
$ws.Range("H86").Value = 2443842.5
$ws.Range("I86").Value = 3227947.5
$ws.Range("K86").Value = 3227947.5
$ws.Range("M86").Value = -3226824.5
$ws.Range("H89").Value = 2443842.5
$ws.Range("I89").Value = 3227947.5
$ws.Range("K89").Value = 16139737.5
$ws.Range("M89").Value = -16134121.5
$ws.Range("H94").Value = 2176956.5
$ws.Range("I94").Value = 2326569.8
$ws.Range("K94").Value = 2326569.8
$ws.Range("M94").Value = -2326118.8
$ws.Range("H105").Value = 3677677.2
$ws.Range("I105").Value = 3907344.5
$ws.Range("K105").Value = 3907344.5
$ws.Range("M105").Value = -3905597.5
$ws.Range("H109").Value = 79995
$ws.Range("J109").Value = 79995
$ws.Range("L109").Value = 79995
$ws.Range("N109").Value = -82769

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22318.34
$ws.Range("I31").Value = 5683.4165
$ws.Range("J31").Value = 36085.17
$ws.Range("K31").Value = 5683.4165
$ws.Range("L31").Value = 36085.17
$ws.Range("M31").Value = -5388.4165
$ws.Range("N31").Value = -36675.17
$ws.Range("H34").Value = 22318.34
$ws.Range("I34").Value = 5683.4165
$ws.Range("J34").Value = 36085.17
$ws.Range("K34").Value = 5683.4165
$ws.Range("L34").Value = 36085.17
$ws.Range("M34").Value = -5481.4165
$ws.Range("N34").Value = -36489.17
$ws.Range("H58").Value = 8802.883
$ws.Range("I58").Value = 10423.091
$ws.Range("J58").Value = 5832.5
$ws.Range("K58").Value = 10423.091
$ws.Range("L58").Value = 5832.5
$ws.Range("M58").Value = -10220.091
$ws.Range("N58").Value = -6238.5
$ws.Range("H99").Value = 3147.4814
$ws.Range("I99").Value = 2560.8235
$ws.Range("J99").Value = 4144.8
$ws.Range("K99").Value = 2560.8235
$ws.Range("L99").Value = 4144.8
$ws.Range("M99").Value = -1062.8235
$ws.Range("N99").Value = -7140.8
$ws.Range("H126").Value = 3147.4814
$ws.Range("I126").Value = 2560.8235
$ws.Range("J126").Value = 4144.8
$ws.Range("K126").Value = 7682.470499999999
$ws.Range("L126").Value = 12434.4
$ws.Range("M126").Value = -5212.470499999999
$ws.Range("N126").Value = -17374.4
$ws.Range("H136").Value = 8802.883
$ws.Range("I136").Value = 10423.091
$ws.Range("J136").Value = 5832.5
$ws.Range("K136").Value = 31269.273
$ws.Range("L136").Value = 17497.5
$ws.Range("M136").Value = -28719.273
$ws.Range("N136").Value = -22597.5
$ws.Range("H141").Value = 228205
$ws.Range("J141").Value = 249995.9
$ws.Range("L141").Value = 249995.9
$ws.Range("N141").Value = -260355.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1923
$ws.Range("I68").Value = 774
$ws.Range("J68").Value = 2251.2856
$ws.Range("K68").Value = 2322
$ws.Range("L68").Value = 6753.8568
$ws.Range("M68").Value = -1511
$ws.Range("N68").Value = -8375.856800000001
$ws.Range("H71").Value = 1923
$ws.Range("I71").Value = 774
$ws.Range("J71").Value = 2251.2856
$ws.Range("K71").Value = 6966
$ws.Range("L71").Value = 20261.5704
$ws.Range("M71").Value = -2910
$ws.Range("N71").Value = -28373.5704

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1253663.6
$ws.Range("I97").Value = 1253663.6
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1253663.6
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1253167.6
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 8417.027
$ws.Range("I132").Value = 6227.276
$ws.Range("K132").Value = 18681.828
$ws.Range("M132").Value = -16151.828

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3343666.2
$ws.Range("I2").Value = 5006750
$ws.Range("J2").Value = 17499.5
$ws.Range("K2").Value = 5006750
$ws.Range("L2").Value = 17499.5
$ws.Range("M2").Value = -5006638
$ws.Range("N2").Value = -17723.5
$ws.Range("H7").Value = 5013.1113
$ws.Range("I7").Value = 2805.611
$ws.Range("K7").Value = 2805.611
$ws.Range("M7").Value = -2693.611
$ws.Range("H40").Value = 7193.615
$ws.Range("I40").Value = 4846.5557
$ws.Range("K40").Value = 4846.5557
$ws.Range("M40").Value = -4710.5557
$ws.Range("H61").Value = 12346124
$ws.Range("I61").Value = 13889264
$ws.Range("K61").Value = 13889264
$ws.Range("M61").Value = -13889062
$ws.Range("H93").Value = 27796522
$ws.Range("I93").Value = 33335828
$ws.Range("K93").Value = 33335828
$ws.Range("M93").Value = -33334580
$ws.Range("H113").Value = 12346124
$ws.Range("I113").Value = 13889264
$ws.Range("K113").Value = 13889264
$ws.Range("M113").Value = -13887094
$ws.Range("H122").Value = 7497.0713
$ws.Range("I122").Value = 4399
$ws.Range("J122").Value = 9218.223
$ws.Range("K122").Value = 13197
$ws.Range("L122").Value = 27654.669
$ws.Range("M122").Value = -10747
$ws.Range("N122").Value = -32554.669
$ws.Range("H126").Value = 5013.1113
$ws.Range("I126").Value = 2805.611
$ws.Range("K126").Value = 8416.832999999999
$ws.Range("M126").Value = -5946.832999999999
$ws.Range("H132").Value = 10761.14
$ws.Range("J132").Value = 5744.8184
$ws.Range("L132").Value = 17234.4552
$ws.Range("N132").Value = -22294.4552
$ws.Range("H136").Value = 52243
$ws.Range("I136").Value = 114151.445
$ws.Range("J136").Value = 5811.6665
$ws.Range("K136").Value = 342454.335
$ws.Range("L136").Value = 17434.9995
$ws.Range("M136").Value = -339904.335
$ws.Range("N136").Value = -22534.9995
$ws.Range("H139").Value = 69027.86
$ws.Range("J139").Value = 69032.5
$ws.Range("L139").Value = 69032.5
$ws.Range("N139").Value = -79312.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3470
$ws.Range("I96").Value = 3278.6
$ws.Range("J96").Value = 3948.5
$ws.Range("K96").Value = 3278.6
$ws.Range("L96").Value = 3948.5
$ws.Range("M96").Value = -1905.6
$ws.Range("N96").Value = -6694.5
$ws.Range("H113").Value = 2098.182
$ws.Range("I113").Value = 1495
$ws.Range("J113").Value = 2158.5
$ws.Range("K113").Value = 4485
$ws.Range("L113").Value = 6475.5
$ws.Range("M113").Value = -2315
$ws.Range("N113").Value = -10815.5
$ws.Range("H126").Value = 3044.6
$ws.Range("I126").Value = 3142.4285
$ws.Range("J126").Value = 2816.3333
$ws.Range("K126").Value = 9427.2855
$ws.Range("L126").Value = 8448.999899999999
$ws.Range("M126").Value = -6957.2855
$ws.Range("N126").Value = -13388.9999
$ws.Range("H132").Value = 18728880
$ws.Range("I132").Value = 22735548
$ws.Range("J132").Value = 1099545.6
$ws.Range("K132").Value = 68206644
$ws.Range("L132").Value = 3298636.8
$ws.Range("M132").Value = -68204114
$ws.Range("H136").Value = 5653.657
$ws.Range("I136").Value = 6987.76
$ws.Range("J136").Value = 2318.4
$ws.Range("K136").Value = 20963.28
$ws.Range("L136").Value = 6955.200000000001
$ws.Range("M136").Value = -18413.28
$ws.Range("N136").Value = -12055.2

